# Insert a new weekly price-report row above row 416 (Feria Lagunitas de
# Puerto Montt - Pepino ensalada), pushing the existing rows 416..510 down
# to 417..511 and extending the used range to A1:R511.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(416).Insert()

$ws.Range("A416").Value = 4
$ws.Range("B416").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C416").Value = "Los Lagos"
$ws.Range("D416").Value = 45204
$ws.Range("E416").Value = 10
$ws.Range("F416").Value = 100112043
$ws.Range("G416").Value = "Pepino ensalada"
$ws.Range("H416").Value = "Sin especificar"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 200
$ws.Range("K416").Value = 21000
$ws.Range("L416").Value = 21000
$ws.Range("M416").Value = 21000
$ws.Range("N416").Value = "$/caja 60 unidades"
$ws.Range("O416").Value = "Región de Arica y Parinacota"
$ws.Range("P416").Value = 350
$ws.Range("Q416").Value = 60
$ws.Range("R416").Value = "Hortaliza"
